$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "div_integrationIcons_internalRoleTreeitemName"
$ws.Range("B1").Value = "div_integrationTools_class"
$ws.Range("C1").Value = "div_integrationTools_class_1"
$ws.Range("D1").Value = "div_menuBackdrop_class"
$ws.Range("E1").Value = "div_menuBackdrop_class_1"
$ws.Range("F1").Value = "input_Name"
$ws.Range("G1").Value = "p_sessionInfo_class"
$ws.Range("H1").Value = "p_sessionInfo_class_1"

# Copy the existing header style (bold "Pandas" style already applied to A1) across the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# --- Data row (row 2) ---
$ws.Range("A2").Value = "Codeless Automation Tool"
$ws.Range("B2").Value = "css-151cthe\ MuiTreeItem-content\ Mui-expanded\ Mui-focused"
$ws.Range("C2").Value = "css-151cthe\ MuiTreeItem-content\ Mui-focused"
$ws.Range("D2").Value = "MuiBackdrop-root\ MuiBackdrop-invisible\ css-esi9ax"
$ws.Range("E2").Value = "MuiBackdrop-root\ MuiBackdrop-invisible\ css-esi9ax"
$ws.Range("G2").Value = "MuiBox-root\ css-0""]:nth-child(5) [class=""MuiTypography-root\ MuiTypography-body1\ css-1lpm9pj"
$ws.Range("H2").Value = "MuiBox-root\ css-0""]:nth-child(1) [class=""MuiTypography-root\ MuiTypography-body1\ css-1lpm9pj"

# F2 stays an empty cell but must still exist in the sheet data (re-apply the default
# style so the cell is preserved even though it carries no value)
$ws.Range("F2").Style = "Normal"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 46.083333333333336
$ws.Columns.Item(2).ColumnWidth = 60.083333333333336
$ws.Columns.Item(3).ColumnWidth = 46.083333333333336
$ws.Columns.Item(4).ColumnWidth = 52.083333333333336
$ws.Columns.Item(5).ColumnWidth = 52.083333333333336
$ws.Columns.Item(6).ColumnWidth = 11.083333333333334
$ws.Columns.Item(7).ColumnWidth = 95.08333333333333
$ws.Columns.Item(8).ColumnWidth = 95.08333333333333
